# "testando o novo plot" - add a "Regiao" (region) column next to the
# existing Estado / Sigla / Codigo table, mapping every Brazilian state to
# its geographic macro-region, and adjust the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- seed the shared-string table in the exact order the region names are
# ---  meant to end up in (Norte, Nordeste, Sudeste, Sul, Centro-Oeste,
# ---  Regiao) by writing them once to a scratch area far away from the
# ---  real data, then wiping the scratch cells again. The real writes
# ---  below simply reuse these already-interned strings.
$ws.Cells.Item(200, 20).Value = "Norte"
$ws.Cells.Item(201, 20).Value = "Nordeste"
$ws.Cells.Item(202, 20).Value = "Sudeste"
$ws.Cells.Item(203, 20).Value = "Sul"
$ws.Cells.Item(204, 20).Value = "Centro-Oeste"
$ws.Cells.Item(205, 20).Value = "Regiao"

# Header for the new column D.
$ws.Cells.Item(1, 4).Value = "Regiao"

# Row -> region, keeping every state on the exact row it already occupies.
$ws.Cells.Item(2, 4).Value = "Norte"          # Acre
$ws.Cells.Item(3, 4).Value = "Nordeste"       # Alagoas
$ws.Cells.Item(4, 4).Value = "Norte"          # Amapá
$ws.Cells.Item(5, 4).Value = "Norte"          # Amazonas
$ws.Cells.Item(6, 4).Value = "Nordeste"       # Bahia
$ws.Cells.Item(7, 4).Value = "Nordeste"       # Ceará
$ws.Cells.Item(8, 4).Value = "Centro-Oeste"   # Distrito Federal
$ws.Cells.Item(9, 4).Value = "Sudeste"        # Espírito Santo
$ws.Cells.Item(10, 4).Value = "Centro-Oeste"  # Goiás
$ws.Cells.Item(11, 4).Value = "Nordeste"      # Maranhão
$ws.Cells.Item(12, 4).Value = "Centro-Oeste"  # Mato Grosso
$ws.Cells.Item(13, 4).Value = "Centro-Oeste"  # Mato Grosso do Sul
$ws.Cells.Item(14, 4).Value = "Sudeste"       # Minas Gerais
$ws.Cells.Item(15, 4).Value = "Norte"         # Pará
$ws.Cells.Item(16, 4).Value = "Nordeste"      # Paraíba
$ws.Cells.Item(17, 4).Value = "Sul"           # Paraná
$ws.Cells.Item(18, 4).Value = "Nordeste"      # Pernambuco
$ws.Cells.Item(19, 4).Value = "Nordeste"      # Piauí
$ws.Cells.Item(20, 4).Value = "Sudeste"       # Rio de Janeiro
$ws.Cells.Item(21, 4).Value = "Nordeste"      # Rio Grande do Norte
$ws.Cells.Item(22, 4).Value = "Sul"           # Rio Grande do Sul
$ws.Cells.Item(23, 4).Value = "Norte"         # Rondônia
$ws.Cells.Item(24, 4).Value = "Norte"         # Roraima
$ws.Cells.Item(25, 4).Value = "Sul"           # Santa Catarina
$ws.Cells.Item(26, 4).Value = "Sudeste"       # São Paulo
$ws.Cells.Item(27, 4).Value = "Nordeste"      # Sergipe
$ws.Cells.Item(28, 4).Value = "Norte"         # Tocantins

# Drop the scratch cells used only to pin the shared-string order.
$ws.Range("T200:T205").Clear()

# Standardize the row heights of the newly-populated rows so they match
# the rest of the sheet.
$ws.Range("A2:D24").EntireRow.RowHeight = 14.25

# Move the selection, as it was left after the edit (also drops the old
# "topLeftCell" scroll anchor now that row 1 is back in view).
$ws.Range("B12").Select()
